$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 303.25
$ws.Range("I2").Value = 322.63635
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 322.63635
$ws.Range("L2").Value = 90
$ws.Range("M2").Value = -209.63635
$ws.Range("N2").Value = -316
$ws.Range("H70").Value = 4818
$ws.Range("I70").Value = 2999.6667
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 8999.000100000001
$ws.Range("L70").Value = 21000
$ws.Range("M70").Value = -8729.000100000001
$ws.Range("N70").Value = -21540
$ws.Range("H73").Value = 4818
$ws.Range("I73").Value = 2999.6667
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 8999.000100000001
$ws.Range("L73").Value = 21000
$ws.Range("M73").Value = -8063.000100000001
$ws.Range("N73").Value = -22872
$ws.Range("H74").Value = 8638.223
$ws.Range("I74").Value = 5998
$ws.Range("K74").Value = 5998
$ws.Range("M74").Value = -5062
$ws.Range("H77").Value = 8638.223
$ws.Range("I77").Value = 5998
$ws.Range("K77").Value = 29990
$ws.Range("M77").Value = -25310
$ws.Range("H80").Value = 2249.75
$ws.Range("I80").Value = 2800
$ws.Range("K80").Value = 8400
$ws.Range("M80").Value = -7402
$ws.Range("H83").Value = 2249.75
$ws.Range("I83").Value = 2800
$ws.Range("K83").Value = 25200
$ws.Range("M83").Value = -20208
$ws.Range("H88").Value = 596.1875
$ws.Range("I88").Value = 439
$ws.Range("J88").Value = 648.5833
$ws.Range("K88").Value = 439
$ws.Range("L88").Value = 648.5833
$ws.Range("M88").Value = -33
$ws.Range("N88").Value = -1460.5833
$ws.Range("H91").Value = 596.1875
$ws.Range("I91").Value = 439
$ws.Range("J91").Value = 648.5833
$ws.Range("K91").Value = 439
$ws.Range("L91").Value = 648.5833
$ws.Range("M91").Value = 965
$ws.Range("N91").Value = -3456.5833
$ws.Range("H94").Value = 5124.1665
$ws.Range("I94").Value = 5124.1665
$ws.Range("K94").Value = 5124.1665
$ws.Range("M94").Value = -4673.1665
$ws.Range("H100").Value = 2440.6
$ws.Range("I100").Value = 1925.1666
$ws.Range("K100").Value = 1925.1666
$ws.Range("M100").Value = -1384.1666
$ws.Range("H129").Value = 2324.9167
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("H137").Value = 4550026
$ws.Range("J137").Value = 6799.7036
$ws.Range("L137").Value = 20399.1108
$ws.Range("N137").Value = -25499.1108
$ws.Range("H138").Value = 5480
$ws.Range("J138").Value = 6100
$ws.Range("L138").Value = 18300
$ws.Range("N138").Value = -28580
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4606.3335
$ws.Range("I61").Value = 2642.1428
$ws.Range("K61").Value = 2642.1428
$ws.Range("M61").Value = -2430.1428
$ws.Range("H88").Value = 3482.3333
$ws.Range("I88").Value = 1799.6666
$ws.Range("K88").Value = 1799.6666
$ws.Range("M88").Value = -1393.6666
$ws.Range("H91").Value = 3482.3333
$ws.Range("I91").Value = 1799.6666
$ws.Range("K91").Value = 1799.6666
$ws.Range("M91").Value = -395.6666
$ws.Range("H97").Value = 2058146.2
$ws.Range("I97").Value = 2179155
$ws.Range("K97").Value = 2179155
$ws.Range("M97").Value = -2178659
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").ClearContents()
$ws.Range("N109").Value = 0
$ws.Range("H136").Value = 4606.3335
$ws.Range("I136").Value = 2642.1428
$ws.Range("K136").Value = 7926.428400000001
$ws.Range("M136").Value = -5376.428400000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2573.353
$ws.Range("J86").Value = 2702.5
$ws.Range("L86").Value = 2702.5
$ws.Range("N86").Value = -4948.5
$ws.Range("H89").Value = 2573.353
$ws.Range("J89").Value = 2702.5
$ws.Range("L89").Value = 13512.5
$ws.Range("N89").Value = -24744.5
$ws.Range("H94").Value = 909.5
$ws.Range("I94").Value = 788.5
$ws.Range("J94").Value = 1070.8334
$ws.Range("K94").Value = 788.5
$ws.Range("L94").Value = 1070.8334
$ws.Range("M94").Value = -337.5
$ws.Range("N94").Value = -1972.8334
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H105").Value = 90912080
$ws.Range("I105").Value = 142859490
$ws.Range("K105").Value = 142859490
$ws.Range("M105").Value = -142857743
$ws.Range("H108").Value = 169994.5
$ws.Range("J108").Value = 169994.5
$ws.Range("L108").Value = 169994.5
$ws.Range("N108").Value = -177674.5
$ws.Range("H134").Value = 5632.852
$ws.Range("J134").Value = 8030.933
$ws.Range("L134").Value = 24092.799
$ws.Range("N134").Value = -29162.799
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3129389.2
$ws.Range("J31").Value = 4767647
$ws.Range("L31").Value = 4767647
$ws.Range("N31").Value = -4768237
$ws.Range("H34").Value = 3129389.2
$ws.Range("J34").Value = 4767647
$ws.Range("L34").Value = 4767647
$ws.Range("N34").Value = -4768051
$ws.Range("H58").Value = 4150.8
$ws.Range("J58").Value = 5593.7144
$ws.Range("L58").Value = 5593.7144
$ws.Range("N58").Value = -5999.7144
$ws.Range("H94").Value = 4630.5713
$ws.Range("J94").Value = 4630.5713
$ws.Range("L94").Value = 4630.5713
$ws.Range("N94").Value = -5532.5713
$ws.Range("H105").Value = 5465.1816
$ws.Range("I105").Value = 5861.7
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 5861.7
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = -4114.7
$ws.Range("N105").Value = -4994
$ws.Range("H136").Value = 4150.8
$ws.Range("J136").Value = 5593.7144
$ws.Range("L136").Value = 16781.1432
$ws.Range("N136").Value = -21881.1432
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2776.3845
$ws.Range("I2").Value = 79.25
$ws.Range("J2").Value = 5088.2144
$ws.Range("K2").Value = 475.5
$ws.Range("L2").Value = 30529.2864
$ws.Range("M2").Value = -362.5
$ws.Range("N2").Value = -30755.2864
$ws.Range("H107").Value = 4560.9
$ws.Range("I107").Value = 1802.5
$ws.Range("J107").Value = 5250.5
$ws.Range("K107").Value = 5407.5
$ws.Range("L107").Value = 15751.5
$ws.Range("M107").Value = -3487.5
$ws.Range("N107").Value = -19591.5
$ws.Range("H122").Value = 17882.8
$ws.Range("J122").Value = 4996.6665
$ws.Range("L122").Value = 44969.9985
$ws.Range("N122").Value = -49869.9985
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1601.8235
$ws.Range("I97").Value = 1587.9286
$ws.Range("K97").Value = 1587.9286
$ws.Range("M97").Value = -1091.9286
$ws.Range("H126").Value = 3986
$ws.Range("J126").Value = 5126.8887
$ws.Range("L126").Value = 15380.6661
$ws.Range("N126").Value = -20320.6661
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 4350
$ws.Range("I4").Value = 3700
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 3700
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = -3587
$ws.Range("N4").Value = -5226
$ws.Range("H28").Value = 4350
$ws.Range("I28").Value = 3700
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 3700
$ws.Range("L28").Value = 5000
$ws.Range("M28").Value = -3468
$ws.Range("N28").Value = -5464
$ws.Range("H37").Value = 4350
$ws.Range("I37").Value = 3700
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 3700
$ws.Range("L37").Value = 5000
$ws.Range("M37").Value = -3593
$ws.Range("N37").Value = -5214
$ws.Range("H48").Value = 15000
$ws.Range("I48").Value = 15000
$ws.Range("K48").Value = 15000
$ws.Range("M48").Value = -14339
$ws.Range("H55").Value = 783.6
$ws.Range("I55").Value = 589.5789
$ws.Range("K55").Value = 589.5789
$ws.Range("M55").Value = -416.5789
$ws.Range("H68").Value = 2857.1428
$ws.Range("I68").Value = 3000.2
$ws.Range("K68").Value = 3000.2
$ws.Range("M68").Value = -2251.2
$ws.Range("H71").Value = 2857.1428
$ws.Range("I71").Value = 3000.2
$ws.Range("K71").Value = 15001
$ws.Range("M71").Value = -11257
$ws.Range("H93").Value = 3427.6875
$ws.Range("I93").Value = 7473.8335
$ws.Range("K93").Value = 7473.8335
$ws.Range("M93").Value = -6225.8335
$ws.Range("H100").Value = 14710323
$ws.Range("I100").Value = 35718000
$ws.Range("J100").Value = 4949.3
$ws.Range("K100").Value = 35718000
$ws.Range("L100").Value = 4949.3
$ws.Range("M100").Value = -35717459
$ws.Range("N100").Value = -6031.3
$ws.Range("H108").Value = 105263
$ws.Range("J108").Value = 105263
$ws.Range("L108").Value = 105263
$ws.Range("N108").Value = -112943
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1869.75
$ws.Range("I96").Value = 1344.5
$ws.Range("K96").Value = 1344.5
$ws.Range("M96").Value = 28.5
$ws.Range("H111").Value = 77000
$ws.Range("J111").Value = 77000
$ws.Range("L111").Value = 77000
$ws.Range("N111").Value = -85180
$ws.Range("H132").Value = 4171.016
$ws.Range("I132").Value = 2338.9546
$ws.Range("K132").Value = 4486.8638
$ws.Range("M132").Value = -4486.8638
